# KP-11725 D: Extension of questionnaire's translation files
# Insert a new "Variable" column right after the "Entity Id" column (A),
# shifting the existing Type/Index/Original/Translation columns one to the
# right, and fill the new column with the constant variable name "c1" for
# every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns B:E -> C:F by inserting a new blank column at B.
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "Variable"

# Fill the new column with the variable name for each translation row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "c1"
}

# Size the new "Variable" column to fit its short content, same as Excel
# does automatically for a freshly typed narrow column.
$ws.Columns.Item(2).ColumnWidth = 6.83

$ws.Range("E6").Select()
